$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was inserted as row 364, pushing the
# existing rows 364-481 down to 365-482 (dimension grows from R481 to R482).
$ws.Rows(364).Insert()

# After the insert, row 365 still holds the data that used to live in the
# old row 364, and row 366 holds the data that used to live in the old row
# 365. The new row 364 reuses most of the old-row-364 values, but gets a
# fresh date (the day after the latest date in the sheet) and picks up the
# "Volumen" (J) / "Origen" (O) values that used to belong to old row 365.
$ws.Range("A364").Value = $ws.Range("A365").Value2
$ws.Range("B364").Value = $ws.Range("B365").Value2
$ws.Range("C364").Value = $ws.Range("C365").Value2
$ws.Range("D364").Value = 44985
$ws.Range("E364").Value = $ws.Range("E365").Value2
$ws.Range("F364").Value = $ws.Range("F365").Value2
$ws.Range("G364").Value = $ws.Range("G365").Value2
$ws.Range("H364").Value = $ws.Range("H365").Value2
$ws.Range("I364").Value = $ws.Range("I365").Value2
$ws.Range("J364").Value = $ws.Range("J366").Value2
$ws.Range("K364").Value = $ws.Range("K365").Value2
$ws.Range("L364").Value = $ws.Range("L365").Value2
$ws.Range("M364").Value = $ws.Range("M365").Value2
$ws.Range("N364").Value = $ws.Range("N365").Value2
$ws.Range("O364").Value = $ws.Range("O366").Value2
$ws.Range("P364").Value = $ws.Range("P365").Value2
$ws.Range("Q364").Value = $ws.Range("Q365").Value2
$ws.Range("R364").Value = $ws.Range("R365").Value2
